# Sniper rounds damage nerf
# Lower the "AP" (H column) multiplier value for several sniper-caliber
# ammo rows on Feuil1. Dependent formulas (J, K, E columns) recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("H2").Value = 1.26
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0.85
$ws.Range("H7").Value = 0.85
$ws.Range("H8").Value = 0.85
$ws.Range("H9").Value = 0.89
$ws.Range("H10").Value = 0.89
$ws.Range("H11").Value = 1.2
$ws.Range("H12").Value = 1.2

$excel.Calculate()

# Update active cell selection to H8 to match the recorded edit location
$ws.Activate()
$ws.Range("H8").Select()
